# GridReynas.xlsx -- "Added crucesEnDiagonal and crucesEnDiagonal methods"
#
# The 8x8 grid in B2:I9 gets filled in with "X,Y" coordinate labels, where X
# is the value of the column header (row 1, B..I = 0..7) and Y is the value
# of the row label (column A, rows 2..9 = 0..7).
#
# The cells are written in two passes (matching the shared-strings insertion
# order recorded in the target workbook):
#   1) crucesEnDiagonal: the diagonal (X == Y) is filled first.
#   2) cruces: the rest of the grid is filled column by column (X), then row
#      by row (Y) within each column, skipping the diagonal cell already set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B", "C", "D", "E", "F", "G", "H", "I")

# --- crucesEnDiagonal: fill the diagonal (X == Y) first ---
for ($d = 0; $d -le 7; $d++) {
    $cellRef = "$($cols[$d])$($d + 2)"
    $ws.Range($cellRef).Value = "$d,$d"
}

# --- cruces: fill the rest of the grid, column by column (X), then row by
#     row (Y) within each column, skipping the diagonal ---
for ($x = 0; $x -le 7; $x++) {
    for ($y = 0; $y -le 7; $y++) {
        if ($x -eq $y) { continue }
        $cellRef = "$($cols[$x])$($y + 2)"
        $ws.Range($cellRef).Value = "$x,$y"
    }
}

# The active selection moved from L10 to D3.
$ws.Range("D3").Select()

# Column B picked up an explicit (custom) width, splitting it out of the
# B:XFD default-width run. The stored width value is quantized by this
# runtime's column-width engine, so we pick the ColumnWidth input that lands
# closest to the original default width (11.42578125) once quantized.
$ws.Columns("B").ColumnWidth = 10.666666666666666
